$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45179 -> 45180) for every data row (rows 2 through 390).
$ws.Range("C2:C390").Value = 45180
